$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Concepts" (2nd tab) -- extend the Concept/isa table with
# Fruit > Citrus > {Lime, Orange} > Limo hierarchy.
# ---------------------------------------------------------------
$wsConcepts = $wb.Worksheets.Item(2)

$wsConcepts.Range("B2").Value = "Identifier"

$wsConcepts.Range("A3").Value = "Fruit"
$wsConcepts.Range("B3").Value = "Fruit"

$wsConcepts.Range("A4").Value = "Citrus"
$wsConcepts.Range("B4").Value = "Citrus"
$wsConcepts.Range("C4").Value = "Fruit"

$wsConcepts.Range("A5").Value = "Lime"
$wsConcepts.Range("B5").Value = "Lime"
$wsConcepts.Range("C5").Value = "Citrus"

$wsConcepts.Range("A6").Value = "Orange"
$wsConcepts.Range("B6").Value = "Orange"
$wsConcepts.Range("C6").Value = "Citrus"

$wsConcepts.Range("A7").Value = "Limo"
$wsConcepts.Range("B7").Value = "Limo"
$wsConcepts.Range("C7").Value = "Lime"

$wsConcepts.Range("A8").Value = "Limo"
$wsConcepts.Range("C8").Value = "Orange"

# ---------------------------------------------------------------
# Sheet "Views" (3rd tab) -- [View]/View header swap + new view rows
# ---------------------------------------------------------------
$wsViews = $wb.Worksheets.Item(3)

$wsViews.Range("A1").Value = "[View]"
$wsViews.Range("A2").Value = "View"

$wsViews.Range("A3").Value = "Europa"
$wsViews.Range("B3").Value = "Fruit"

$wsViews.Range("A4").Value = "Spanje"
$wsViews.Range("B4").Value = "Citrus"

$wsViews.Range("E5").Value = "Andalusie"
$wsViews.Range("F5").Value = "Lime"

$wsViews.Range("A6").Value = "Israel"
$wsViews.Range("B6").Value = "Orange"

$wsViews.Range("A7").Value = "Mars"
$wsViews.Range("B7").Value = "Limo"

$wsViews.Range("E5:F5").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Atoms" (1st tab) -- extend the population table; new rows
# reuse the "@" text format already applied to the first few rows.
# ---------------------------------------------------------------
$wsAtoms = $wb.Worksheets.Item(1)

$wsAtoms.Range("C1").Value = "pop"
$wsAtoms.Range("B2").Value = "Representation"

# rows 3-8 keep the existing Text ("@") number format
$wsAtoms.Range("A3:B8").NumberFormat = "@"

$wsAtoms.Range("A3").Value = "een"
$wsAtoms.Range("B3").Value = "een"
$wsAtoms.Range("C3").Value = "Fruit"

$wsAtoms.Range("A4").Value = "twee"
$wsAtoms.Range("B4").Value = "twee"
$wsAtoms.Range("C4").Value = "Citrus"

$wsAtoms.Range("A5").Value = "twee"
$wsAtoms.Range("B5").ClearContents()
$wsAtoms.Range("C5").Value = "Fruit"

$wsAtoms.Range("A6").Value = "drie"
$wsAtoms.Range("B6").Value = "drie"
$wsAtoms.Range("C6").Value = "Lime"

$wsAtoms.Range("A7").Value = "drie"
$wsAtoms.Range("B7").ClearContents()
$wsAtoms.Range("C7").Value = "Citrus"

$wsAtoms.Range("A8").Value = "drie"
$wsAtoms.Range("B8").ClearContents()
$wsAtoms.Range("C8").Value = "Fruit"

$wsAtoms.Range("A9").Value = "vier"
$wsAtoms.Range("B9").Value = "vier"
$wsAtoms.Range("C9").Value = "Orange"

$wsAtoms.Range("A10").Value = "vier"
$wsAtoms.Range("C10").Value = "Citrus"

$wsAtoms.Range("A11").Value = "vier"
$wsAtoms.Range("C11").Value = "Fruit"

$wsAtoms.Range("A12").Value = "zes"
$wsAtoms.Range("B12").Value = "zes"
$wsAtoms.Range("C12").Value = "Limo"

$wsAtoms.Range("A13").Value = "zes"
$wsAtoms.Range("C13").Value = "Lime"

$wsAtoms.Range("A14").Value = "zes"
$wsAtoms.Range("C14").Value = "Orange"

$wsAtoms.Range("A15").Value = "zes"
$wsAtoms.Range("C15").Value = "Citrus"

$wsAtoms.Range("A16").Value = "zes"
$wsAtoms.Range("C16").Value = "Fruit"

$wsAtoms.Range("B3").Select() | Out-Null
